# Apply the numeric/text cell updates described by the commit diff.
# Workbook already loaded; operate on $excel.ActiveWorkbook directly.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('展览')
$ws1.Range('F2').Value = 6819
$ws1.Range('F3').Value = 834
$ws1.Range('F4').Value = 151
$ws1.Range('G4').Value = 70
$ws1.Range('F5').Value = 18
$ws1.Range('F6').Value = 762
$ws1.Range('F7').Value = 762
$ws1.Range('F9').Value = 225
$ws1.Range('F11').Value = 1121
$ws1.Range('F12').Value = 879
$ws1.Range('F13').Value = 15
$ws1.Range('F14').Value = 708
$ws1.Range('F16').Value = 1373
$ws1.Range('F19').Value = 1552
$ws1.Range('F20').Value = 9
$ws1.Range('F21').Value = 585
$ws1.Range('F24').Value = 369
$ws1.Range('F25').Value = 1071
$ws1.Range('F27').Value = 734
$ws1.Range('F28').Value = 567
$ws1.Range('F29').Value = 475
$ws1.Range('F30').Value = 461
$ws1.Range('F32').Value = 1013
$ws1.Range('F33').Value = 1139
$ws1.Range('F34').Value = 281
$ws1.Range('F35').Value = 2380
$ws1.Range('F37').Value = 1298
$ws1.Range('F38').Value = 450
$ws1.Range('F40').Value = 3905

$ws2 = $wb.Worksheets.Item('演出')
$ws2.Range('F2').Value = 752
$ws2.Range('G2').Value = '不可售'
$ws2.Range('F4').Value = 27
$ws2.Range('F5').Value = 1034
$ws2.Range('F7').Value = 172
$ws2.Range('F12').Value = 7
$ws2.Range('F13').Value = 396
$ws2.Range('F14').Value = 341
$ws2.Range('F20').Value = 3
$ws2.Range('F21').Value = 246
$ws2.Range('F23').Value = 115
$ws2.Range('F25').Value = 231
$ws2.Range('F26').Value = 41

$ws3 = $wb.Worksheets.Item('本地生活')
$ws3.Range('F5').Value = 1651

$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Range('F4').Value = 1651
$ws4.Range('F8').Value = 6819
$ws4.Range('F9').Value = 834
$ws4.Range('C10').Value = '上海·次元裂缝-X 新春二次元DJ派对'
$ws4.Range('D10').Value = '海潮路133号B1 JUMP工坊'
$ws4.Range('E10').Value = '2024.03.10 14:00-03.10 19:00'
$ws4.Range('F10').Value = 151
$ws4.Range('G10').Value = 70
$ws4.Range('H10').Value = 'https://show.bilibili.com/platform/detail.html?id=81959'
$ws4.Range('I10').Value = '//i1.hdslb.com/bfs/openplatform/202402/MaO7WWLr1708482746780.jpeg'
$ws4.Range('B11').NumberFormat = "@"
$ws4.Range('B11').Value = '2024-03-15'
$ws4.Range('C11').Value = '上海·坏孩纸物语の第35届动漫节之全民宅舞'
$ws4.Range('D11').Value = '泸定路388号 桃源π商业广场'
$ws4.Range('E11').Value = '2024.03.15 11:30-03.16 16:00'
$ws4.Range('F11').Value = 18
$ws4.Range('G11').Value = 20
$ws4.Range('H11').Value = 'https://show.bilibili.com/platform/detail.html?id=82477'
$ws4.Range('I11').Value = '//i2.hdslb.com/bfs/openplatform/202403/WHufQNn91709782559844.png'
$ws4.Range('B12').NumberFormat = "@"
$ws4.Range('B12').Value = '2024-03-16'
$ws4.Range('C12').Value = '上海·Look Look动漫嘉年华'
$ws4.Range('D12').Value = '龙吴路4800号2号门 有只怪兽片场'
$ws4.Range('E12').Value = '2024.03.16 10:00-03.17 17:30'
$ws4.Range('F12').Value = 762
$ws4.Range('G12').Value = 52.2
$ws4.Range('H12').Value = 'https://show.bilibili.com/platform/detail.html?id=81804'
$ws4.Range('I12').Value = '//i1.hdslb.com/bfs/openplatform/202403/om0OCpxy1709287210276.jpeg'
$ws4.Range('F13').Value = 762
$ws4.Range('C14').Value = '上海·SISP动漫游戏嘉年华'
$ws4.Range('D14').Value = '年家浜路518号 周浦万达广场'
$ws4.Range('E14').Value = '2024.03.16 13:00-03.17 19:00'
$ws4.Range('F14').Value = 225
$ws4.Range('G14').Value = 48
$ws4.Range('H14').Value = 'https://show.bilibili.com/platform/detail.html?id=80339'
$ws4.Range('I14').Value = '//i0.hdslb.com/bfs/openplatform/202312/a8iuOufB1703832570508.jpeg'
$ws4.Range('C15').Value = '上海·次元裂缝·X 二次元DJ派对'
$ws4.Range('D15').Value = '海潮路133号B1 JUMP工坊'
$ws4.Range('E15').Value = '2024.03.16 14:00-03.16 19:00'
$ws4.Range('F15').Value = 29
$ws4.Range('G15').Value = 60
$ws4.Range('H15').Value = 'https://show.bilibili.com/platform/detail.html?id=82359'
$ws4.Range('I15').Value = '//i0.hdslb.com/bfs/openplatform/202403/bLFFO59L1709629243557.jpeg'
$ws4.Range('C16').Value = '上海·第九届ACBC动漫盛典'
$ws4.Range('D16').Value = '漕溪北路339号百脑汇4楼 百脑汇'
$ws4.Range('E16').Value = '2024.03.16 10:00-03.17 18:00'
$ws4.Range('F16').Value = 1121
$ws4.Range('G16').Value = 48.8
$ws4.Range('H16').Value = 'https://show.bilibili.com/platform/detail.html?id=82135'
$ws4.Range('I16').Value = '//i0.hdslb.com/bfs/openplatform/202402/bXTNHlWS1709175765881.jpeg'
$ws4.Range('C17').Value = '上海·第五人格ONLY'
$ws4.Range('D17').Value = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws4.Range('E17').Value = '2024.03.16 10:00-03.16 17:00'
$ws4.Range('F17').Value = 879
$ws4.Range('G17').Value = 60
$ws4.Range('H17').Value = 'https://show.bilibili.com/platform/detail.html?id=81533'
$ws4.Range('I17').Value = '//i1.hdslb.com/bfs/openplatform/202401/sOMO7Bjc1706604737277.png'
$ws4.Range('B18').NumberFormat = "@"
$ws4.Range('B18').Value = '2024-03-17'
$ws4.Range('C18').Value = '上海 ·《疯狂动物城》动漫视听音乐会'
$ws4.Range('D18').Value = '牛庄路704号 中国大戏院'
$ws4.Range('E18').Value = '2024.03.17 15:30-03.17 17:00'
$ws4.Range('F18').Value = 27
$ws4.Range('G18').Value = 80
$ws4.Range('H18').Value = 'https://show.bilibili.com/platform/detail.html?id=81112'
$ws4.Range('I18').Value = '//i2.hdslb.com/bfs/openplatform/202401/Wg8b6SRn1705651166088.png'
$ws4.Range('F19').Value = 708
$ws4.Range('F20').Value = 172
$ws4.Range('F21').Value = 172
$ws4.Range('F24').Value = 1373
$ws4.Range('F27').Value = 1552
$ws4.Range('F28').Value = 9
$ws4.Range('F29').Value = 585
$ws4.Range('F30').Value = 7
$ws4.Range('F31').Value = 341
$ws4.Range('F32').Value = 369
$ws4.Range('F33').Value = 1071
$ws4.Range('F35').Value = 734
$ws4.Range('F36').Value = 567
$ws4.Range('F37').Value = 475
$ws4.Range('F38').Value = 461
$ws4.Range('F42').Value = 1013
$ws4.Range('F43').Value = 1139
$ws4.Range('F44').Value = 281
$ws4.Range('F45').Value = 2380
$ws4.Range('F46').Value = 231
$ws4.Range('F49').Value = 1298
$ws4.Range('F50').Value = 450
$ws4.Range('F51').Value = 3905

